$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 51.93629233333333
$ws.Range("H2").Value = 155.808877
$ws.Range("I2").Value = 0.7704232182162135
$ws.Range("J2").Value = 0.7704232182162134
$ws.Range("M2").Value = 8.226724333333333
$ws.Range("N2").Value = 24.680173
$ws.Range("O2").Value = 0.06198126651953669
$ws.Range("P2").Value = 0.06198126651953669
$ws.Range("Q2").Value = 427.2655599217468
$ws.Range("R2").Value = 3845.390039295721
$ws.Range("S2").Value = 0.04775180682109831
$ws.Range("T2").Value = 0.0477518068210983
$ws.Range("G3").Value = 51.93629233333333
$ws.Range("H3").Value = 155.808877
$ws.Range("I3").Value = 0.7704232182162135
$ws.Range("J3").Value = 0.7704232182162134
$ws.Range("O3").Value = 0.6623065855236785
$ws.Range("P3").Value = 0.6623065855236785
$ws.Range("Q3").Value = 4565.585861567354
$ws.Range("R3").Value = 41090.27275410619
$ws.Range("S3").Value = 0.5102563710649443
$ws.Range("T3").Value = 0.5102563710649441
$ws.Range("G4").Value = 51.93629233333333
$ws.Range("H4").Value = 155.808877
$ws.Range("I4").Value = 0.7704232182162135
$ws.Range("J4").Value = 0.7704232182162134
$ws.Range("M4").Value = 36.43008433333333
$ws.Range("N4").Value = 109.290253
$ws.Range("O4").Value = 0.2744692388979848
$ws.Range("P4").Value = 0.2744692388979848
$ws.Range("Q4").Value = 1892.043509663987
$ws.Range("R4").Value = 17028.39158697588
$ws.Range("S4").Value = 0.2114574743331402
$ws.Range("T4").Value = 0.2114574743331402
$ws.Range("G5").Value = 51.93629233333333
$ws.Range("H5").Value = 155.808877
$ws.Range("I5").Value = 0.7704232182162135
$ws.Range("J5").Value = 0.7704232182162134
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1649703333333333
$ws.Range("N5").Value = 0.494911
$ws.Range("O5").Value = 0.00124290905879997
$ws.Range("P5").Value = 0.00124290905879997
$ws.Range("Q5").Value = 8.567947458327444
$ws.Range("R5").Value = 77.111527124947
$ws.Range("S5").Value = 0.0009575659970307575
$ws.Range("T5").Value = 0.0009575659970307574
$ws.Range("G6").Value = 0.03171066666666666
$ws.Range("H6").Value = 0.09513199999999999
$ws.Range("I6").Value = 0.0004703961867034368
$ws.Range("J6").Value = 0.0004703961867034368
$ws.Range("M6").Value = 8.226724333333333
$ws.Range("N6").Value = 24.680173
$ws.Range("O6").Value = 0.06198126651953669
$ws.Range("P6").Value = 0.06198126651953669
$ws.Range("Q6").Value = 0.2608749130928888
$ws.Range("R6").Value = 2.347874217836
$ws.Range("S6").Value = 0.00002915575141783945
$ws.Range("T6").Value = 0.00002915575141783945
$ws.Range("G7").Value = 0.03171066666666666
$ws.Range("H7").Value = 0.09513199999999999
$ws.Range("I7").Value = 0.0004703961867034368
$ws.Range("J7").Value = 0.0004703961867034368
$ws.Range("O7").Value = 0.6623065855236785
$ws.Range("P7").Value = 0.6623065855236785
$ws.Range("S7").Value = 0.000311546492258912
$ws.Range("T7").Value = 0.000311546492258912
$ws.Range("G8").Value = 0.03171066666666666
$ws.Range("H8").Value = 0.09513199999999999
$ws.Range("I8").Value = 0.0004703961867034368
$ws.Range("J8").Value = 0.0004703961867034368
$ws.Range("M8").Value = 36.43008433333333
$ws.Range("N8").Value = 109.290253
$ws.Range("O8").Value = 0.2744692388979848
$ws.Range("P8").Value = 0.2744692388979848
$ws.Range("Q8").Value = 1.155222260932889
$ws.Range("R8").Value = 10.397000348396
$ws.Range("S8").Value = 0.0001291092833450066
$ws.Range("T8").Value = 0.0001291092833450066
$ws.Range("G9").Value = 0.03171066666666666
$ws.Range("H9").Value = 0.09513199999999999
$ws.Range("I9").Value = 0.0004703961867034368
$ws.Range("J9").Value = 0.0004703961867034368
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1649703333333333
$ws.Range("N9").Value = 0.494911
$ws.Range("O9").Value = 0.00124290905879997
$ws.Range("P9").Value = 0.00124290905879997
$ws.Range("Q9").Value = 0.005231319250222222
$ws.Range("R9").Value = 0.04708187325199999
$ws.Range("S9").Value = 0.0000005846596816786634
$ws.Range("T9").Value = 0.0000005846596816786634
$ws.Range("G10").Value = 2.162051666666667
$ws.Range("H10").Value = 6.486155
$ws.Range("I10").Value = 0.03207188515291837
$ws.Range("J10").Value = 0.03207188515291837
$ws.Range("M10").Value = 8.226724333333333
$ws.Range("N10").Value = 24.680173
$ws.Range("O10").Value = 0.06198126651953669
$ws.Range("P10").Value = 0.06198126651953669
$ws.Range("Q10").Value = 17.78660305609056
$ws.Range("R10").Value = 160.079427504815
$ws.Range("S10").Value = 0.001987856061447005
$ws.Range("T10").Value = 0.001987856061447005
$ws.Range("G11").Value = 2.162051666666667
$ws.Range("H11").Value = 6.486155
$ws.Range("I11").Value = 0.03207188515291837
$ws.Range("J11").Value = 0.03207188515291837
$ws.Range("O11").Value = 0.6623065855236785
$ws.Range("P11").Value = 0.6623065855236785
$ws.Range("Q11").Value = 190.0604004990961
$ws.Range("R11").Value = 1710.543604491865
$ws.Range("S11").Value = 0.02124142074693692
$ws.Range("T11").Value = 0.02124142074693692
$ws.Range("G12").Value = 2.162051666666667
$ws.Range("H12").Value = 6.486155
$ws.Range("I12").Value = 0.03207188515291837
$ws.Range("J12").Value = 0.03207188515291837
$ws.Range("M12").Value = 36.43008433333333
$ws.Range("N12").Value = 109.290253
$ws.Range("O12").Value = 0.2744692388979848
$ws.Range("P12").Value = 0.2744692388979848
$ws.Range("Q12").Value = 78.76372454969055
$ws.Range("R12").Value = 708.8735209472151
$ws.Range("S12").Value = 0.008802745907945083
$ws.Range("T12").Value = 0.008802745907945083
$ws.Range("G13").Value = 2.162051666666667
$ws.Range("H13").Value = 6.486155
$ws.Range("I13").Value = 0.03207188515291837
$ws.Range("J13").Value = 0.03207188515291837
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1649703333333333
$ws.Range("N13").Value = 0.494911
$ws.Range("O13").Value = 0.00124290905879997
$ws.Range("P13").Value = 0.00124290905879997
$ws.Range("Q13").Value = 0.3566743841338889
$ws.Range("R13").Value = 3.210069457205
$ws.Range("S13").Value = 0.00003986243658935448
$ws.Range("T13").Value = 0.00003986243658935448
$ws.Range("G14").Value = 13.28262333333333
$ws.Range("H14").Value = 39.84787
$ws.Range("I14").Value = 0.1970345004441647
$ws.Range("J14").Value = 0.1970345004441647
$ws.Range("M14").Value = 8.226724333333333
$ws.Range("N14").Value = 24.680173
$ws.Range("O14").Value = 0.06198126651953669
$ws.Range("P14").Value = 0.06198126651953669
$ws.Range("Q14").Value = 109.2724805868344
$ws.Range("R14").Value = 983.45232528151
$ws.Range("S14").Value = 0.01221244788557354
$ws.Range("T14").Value = 0.01221244788557354
$ws.Range("G15").Value = 13.28262333333333
$ws.Range("H15").Value = 39.84787
$ws.Range("I15").Value = 0.1970345004441647
$ws.Range("J15").Value = 0.1970345004441647
$ws.Range("O15").Value = 0.6623065855236785
$ws.Range("P15").Value = 0.6623065855236785
$ws.Range("Q15").Value = 1167.641249898579
$ws.Range("R15").Value = 10508.77124908721
$ws.Range("S15").Value = 0.1304972472195385
$ws.Range("T15").Value = 0.1304972472195385
$ws.Range("G16").Value = 13.28262333333333
$ws.Range("H16").Value = 39.84787
$ws.Range("I16").Value = 0.1970345004441647
$ws.Range("J16").Value = 0.1970345004441647
$ws.Range("M16").Value = 36.43008433333333
$ws.Range("N16").Value = 109.290253
$ws.Range("O16").Value = 0.2744692388979848
$ws.Range("P16").Value = 0.2744692388979848
$ws.Range("Q16").Value = 483.8870882012345
$ws.Range("R16").Value = 4354.98379381111
$ws.Range("S16").Value = 0.05407990937355454
$ws.Range("T16").Value = 0.05407990937355454
$ws.Range("G17").Value = 13.28262333333333
$ws.Range("H17").Value = 39.84787
$ws.Range("I17").Value = 0.1970345004441647
$ws.Range("J17").Value = 0.1970345004441647
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1649703333333333
$ws.Range("N17").Value = 0.494911
$ws.Range("O17").Value = 0.00124290905879997
$ws.Range("P17").Value = 0.00124290905879997
$ws.Range("Q17").Value = 2.191238798841111
$ws.Range("R17").Value = 19.72114918957
$ws.Range("S17").Value = 0.000244895965498179
$ws.Range("T17").Value = 0.000244895965498179
